$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("X2").Value = "Utility (Percent)"
$ws.Range("L3").Value = "20 msec"
$ws.Range("M3").Value = "2734 msec"
$ws.Range("N3").Value = "1275.62 msec"
$ws.Range("O3").Value = "1301 usec"
$ws.Range("P3").Value = "1637.5k usec"
$ws.Range("Q3").Value = "8686.07 usec"
$ws.Range("L4").Value = "5 msec"
$ws.Range("M4").Value = "2464 msec"
$ws.Range("N4").Value = "1272.88 msec"
$ws.Range("O4").Value = "1362 usec"
$ws.Range("P4").Value = "1705.3k usec"
$ws.Range("Q4").Value = "14938.35 usec"
$ws.Range("L5").Value = "647 usec"
$ws.Range("M5").Value = "5327.3k usec"
$ws.Range("N5").Value = "1253619.34 usec"
$ws.Range("O5").Value = "1261 usec"
$ws.Range("P5").Value = "6556.3k usec"
$ws.Range("Q5").Value = "74929.02 usec"
$ws.Range("L6").Value = "5 msec"
$ws.Range("M6").Value = "2881 msec"
$ws.Range("N6").Value = "1274.55 msec"
$ws.Range("O6").Value = "1383 usec"
$ws.Range("P6").Value = "1848.1k usec"
$ws.Range("Q6").Value = "11302.20 usec"
$ws.Range("L7").Value = "11 msec"
$ws.Range("M7").Value = "2911 msec"
$ws.Range("N7").Value = "1255.77 msec"
$ws.Range("O7").Value = "1386 usec"
$ws.Range("P7").Value = "2496.2k usec"
$ws.Range("Q7").Value = "55184.95 usec"
$ws.Range("L8").Value = "325 usec"
$ws.Range("M8").Value = "9296.7k usec"
$ws.Range("N8").Value = "1060430.31 usec"
$ws.Range("O8").Value = "1765 usec"
$ws.Range("P8").Value = "12929k usec"
$ws.Range("Q8").Value = "1297209.56 usec"
$ws.Range("L9").Value = "330 usec"
$ws.Range("M9").Value = "5553.4k usec"
$ws.Range("N9").Value = "1271092.60 usec"
$ws.Range("O9").Value = "1392 usec"
$ws.Range("P9").Value = "8060.7k usec"
$ws.Range("Q9").Value = "88738.80 usec"
$ws.Range("L10").Value = "9 msec"
$ws.Range("M10").Value = "3492 msec"
$ws.Range("N10").Value = "1275.43 msec"
$ws.Range("O10").Value = "1418 usec"
$ws.Range("P10").Value = "2415.0k usec"
$ws.Range("Q10").Value = "8984.37 usec"
$ws.Range("L11").Value = "9 msec"
$ws.Range("M11").Value = "2489 msec"
$ws.Range("N11").Value = "1277.69 msec"
$ws.Range("O11").Value = "1393 usec"
$ws.Range("P11").Value = "620273 usec"
$ws.Range("Q11").Value = "3842.78 usec"
$ws.Range("L12").Value = "6 msec"
$ws.Range("M12").Value = "2500 msec"
$ws.Range("N12").Value = "1272.65 msec"
$ws.Range("O12").Value = "2 msec"
$ws.Range("P12").Value = "514 msec"
$ws.Range("Q12").Value = "15.82 msec"
$ws.Range("L13").Value = "1482 usec"
$ws.Range("M13").Value = "2485.3k usec"
$ws.Range("N13").Value = "1275156.71 usec"
$ws.Range("O13").Value = "1372 usec"
$ws.Range("P13").Value = "688717 usec"
$ws.Range("Q13").Value = "9379.13 usec"
$ws.Range("L14").Value = "2 msec"
$ws.Range("M14").Value = "2449 msec"
$ws.Range("N14").Value = "1270.54 msec"
$ws.Range("O14").Value = "2 msec"
$ws.Range("P14").Value = "1487 msec"
$ws.Range("Q14").Value = "20.04 msec"
$ws.Range("L15").Value = "330 usec"
$ws.Range("M15").Value = "5232.0k usec"
$ws.Range("N15").Value = "851725.78 usec"
$ws.Range("O15").Value = "1567 usec"
$ws.Range("P15").Value = "8542.8k usec"
$ws.Range("Q15").Value = "1357253.92 usec"
$ws.Range("L16").Value = "306 usec"
$ws.Range("M16").Value = "4194.9k usec"
$ws.Range("N16").Value = "1211201.74 usec"
$ws.Range("O16").Value = "1443 usec"
$ws.Range("P16").Value = "5451.8k usec"
$ws.Range("Q16").Value = "193189.34 usec"
$ws.Range("L17").Value = "7 msec"
$ws.Range("M17").Value = "2465 msec"
$ws.Range("N17").Value = "1274.15 msec"
$ws.Range("O17").Value = "1594 usec"
$ws.Range("P17").Value = "1519.0k usec"
$ws.Range("Q17").Value = "12287.02 usec"
$ws.Range("L18").Value = "526 usec"
$ws.Range("M18").Value = "5794.3k usec"
$ws.Range("N18").Value = "1262881.16 usec"
$ws.Range("O18").Value = "1230 usec"
$ws.Range("P18").Value = "6902.5k usec"
$ws.Range("Q18").Value = "57778.46 usec"
$ws.Range("L19").Value = "529 usec"
$ws.Range("M19").Value = "3832.8k usec"
$ws.Range("N19").Value = "1141353.26 usec"
$ws.Range("O19").Value = "1383 usec"
$ws.Range("P19").Value = "5156.5k usec"
$ws.Range("Q19").Value = "323333.12 usec"
$ws.Range("L20").Value = "1090 usec"
$ws.Range("M20").Value = "2536.8k usec"
$ws.Range("N20").Value = "1268200.07 usec"
$ws.Range("O20").Value = "2 msec"
$ws.Range("P20").Value = "1429 msec"
$ws.Range("Q20").Value = "26.22 msec"
$ws.Range("L21").Value = "297 usec"
$ws.Range("M21").Value = "5232.7k usec"
$ws.Range("N21").Value = "687001.26 usec"
$ws.Range("O21").Value = "2 msec"
$ws.Range("P21").Value = "6784 msec"
$ws.Range("Q21").Value = "1844.93 msec"
$ws.Range("L22").Value = "9 msec"
$ws.Range("M22").Value = "2590 msec"
$ws.Range("N22").Value = "1245.47 msec"
$ws.Range("O22").Value = "1262 usec"
$ws.Range("P22").Value = "2159.5k usec"
$ws.Range("Q22").Value = "79003.52 usec"
$ws.Range("L23").Value = "801 usec"
$ws.Range("M23").Value = "2483.8k usec"
$ws.Range("N23").Value = "1271731.74 usec"
$ws.Range("O23").Value = "1393 usec"
$ws.Range("P23").Value = "1582.5k usec"
$ws.Range("Q23").Value = "17702.64 usec"
